# Automatic update of files.
# Column C ("Förändrad") on the active sheet holds the "last changed" date
# stamp for every data row (rows 2-33). Bump it from 45615 (2024-11-19) to
# 45616 (2024-11-20) for each row that currently carries the old stamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45615) {
        $cell.Value = 45616
    }
}
